$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "JALR_o" column (M) header, mirroring the other signal-header cells (row 2)
$ws.Range("M2").Value = "JALR_o"

# Data rows 3-12: JALR_o is active ("1'b1") only for the JALR instruction (row 7),
# "1'b0" everywhere else. Apply the same numeric display format ("0") used by the
# other data cells in these rows so the new cells share style s="1".
$ws.Range("M3").Value = "1'b0"
$ws.Range("M4").Value = "1'b0"
$ws.Range("M5").Value = "1'b0"
$ws.Range("M6").Value = "1'b0"
$ws.Range("M7").Value = "1'b1"
$ws.Range("M8").Value = "1'b0"
$ws.Range("M9").Value = "1'b0"
$ws.Range("M10").Value = "1'b0"
$ws.Range("M11").Value = "1'b0"
$ws.Range("M12").Value = "1'b0"

$ws.Range("M3:M12").NumberFormat = "0"

# Update the saved selection on the sheet
$ws.Range("N7").Select()
